$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I ("I0") and J ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold font + border + centered alignment) from H1
# onto the two new header cells, matching the rest of the header row
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill in the new I0 / IF data columns for every data row (2-71)
$data = @(
    ,@(2, 8, 9)
    ,@(3, 7, 7)
    ,@(4, 8, 8)
    ,@(5, 6, 7)
    ,@(6, 6, 7)
    ,@(7, 9, 9)
    ,@(8, 3, 4)
    ,@(9, 6, 6)
    ,@(10, 6, 7)
    ,@(11, 8, 8)
    ,@(12, 8, 8)
    ,@(13, 9, 9)
    ,@(14, 8, 8)
    ,@(15, 6, 6)
    ,@(16, 6, 6)
    ,@(17, 9, 9)
    ,@(18, 8, 8)
    ,@(19, 6, 7)
    ,@(20, 7, 7)
    ,@(21, 8, 8)
    ,@(22, 7, 7)
    ,@(23, 6, 7)
    ,@(24, 8, 8)
    ,@(25, 5, 5)
    ,@(26, 6, 6)
    ,@(27, 10, 10)
    ,@(28, 9, 9)
    ,@(29, 7, 7)
    ,@(30, 5, 6)
    ,@(31, 7, 7)
    ,@(32, 8, 8)
    ,@(33, 6, 6)
    ,@(34, 8, 9)
    ,@(35, 6, 7)
    ,@(36, 8, 8)
    ,@(37, 8, 8)
    ,@(38, 6, 7)
    ,@(39, 9, 9)
    ,@(40, 4, 4)
    ,@(41, 10, 10)
    ,@(42, 5, 5)
    ,@(43, 7, 7)
    ,@(44, 7, 7)
    ,@(45, 10, 10)
    ,@(46, 3, 3)
    ,@(47, 7, 7)
    ,@(48, 9, 9)
    ,@(49, 3, 4)
    ,@(50, 7, 7)
    ,@(51, 8, 8)
    ,@(52, 9, 10)
    ,@(53, 7, 7)
    ,@(54, 6, 7)
    ,@(55, 6, 6)
    ,@(56, 6, 6)
    ,@(57, 7, 7)
    ,@(58, 9, 9)
    ,@(59, 11, 11)
    ,@(60, 8, 8)
    ,@(61, 10, 10)
    ,@(62, 6, 7)
    ,@(63, 5, 5)
    ,@(64, 8, 8)
    ,@(65, 5, 5)
    ,@(66, 6, 6)
    ,@(67, 8, 8)
    ,@(68, 5, 5)
    ,@(69, 8, 8)
    ,@(70, 5, 5)
    ,@(71, 5, 5)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
